# Scheduled-runner refresh of market-price-derived columns (H:N) across all
# eight job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Source data is
# fetched live market pricing, so only numeric cells move; no formulas,
# headers, or structural changes are involved.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3536.611
$ws.Range("I40").Value = 4376.4
$ws.Range("J40").Value = 3213.6155
$ws.Range("K40").Value = 4376.4
$ws.Range("L40").Value = 3213.6155
$ws.Range("M40").Value = -4201.4
$ws.Range("N40").Value = -3563.6155
$ws.Range("H43").Value = 3413.4546
$ws.Range("I43").Value = 4018.5
$ws.Range("J43").Value = 1800
$ws.Range("K43").Value = 4018.5
$ws.Range("L43").Value = 1800
$ws.Range("M43").Value = -3949.5
$ws.Range("N43").Value = -1938
$ws.Range("H51").Value = 8873.4
$ws.Range("J51").Value = 8591.75
$ws.Range("L51").Value = 8591.75
$ws.Range("N51").Value = -9559.75
$ws.Range("H137").Value = 22786.22
$ws.Range("I137").Value = 26336.23
$ws.Range("K137").Value = 79008.69
$ws.Range("M137").Value = -76458.69

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 95
$ws.Range("I5").Value = 95
$ws.Range("K5").Value = 95
$ws.Range("M5").Value = 17
$ws.Range("H45").Value = 36711.617
$ws.Range("I45").Value = 45523.4
$ws.Range("J45").Value = 7339
$ws.Range("K45").Value = 45523.4
$ws.Range("L45").Value = 7339
$ws.Range("M45").Value = -45146.4
$ws.Range("N45").Value = -8093
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H110").Value = 2689.963
$ws.Range("I110").Value = 1573.381
$ws.Range("J110").Value = 6598
$ws.Range("K110").Value = 1573.381
$ws.Range("L110").Value = 6598
$ws.Range("M110").Value = 471.6189999999999
$ws.Range("N110").Value = -10688
$ws.Range("H122").Value = 4413.4243
$ws.Range("I122").Value = 3919.8076
$ws.Range("K122").Value = 11759.4228
$ws.Range("M122").Value = -9309.4228
$ws.Range("H126").Value = 5999.5
$ws.Range("I126").Value = 5999.5
$ws.Range("K126").Value = 17998.5
$ws.Range("M126").Value = -15528.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 95
$ws.Range("I4").Value = 95
$ws.Range("K4").Value = 95
$ws.Range("M4").Value = 20
$ws.Range("H105").Value = 16252948
$ws.Range("I105").Value = 1002371.2
$ws.Range("K105").Value = 1002371.2
$ws.Range("M105").Value = -1000624.2
$ws.Range("H107").Value = 1207.5454
$ws.Range("I107").Value = 1218.3
$ws.Range("K107").Value = 1218.3
$ws.Range("M107").Value = 701.7
$ws.Range("H113").Value = 6000
$ws.Range("I113").Value = 6000
$ws.Range("K113").Value = 6000
$ws.Range("M113").Value = -3830

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1558.5714
$ws.Range("I22").Value = 1318.5
$ws.Range("J22").Value = 2999
$ws.Range("K22").Value = 1318.5
$ws.Range("L22").Value = 2999
$ws.Range("M22").Value = -968.5
$ws.Range("N22").Value = -3699
$ws.Range("H62").Value = 11120868
$ws.Range("J62").Value = 12518.2
$ws.Range("L62").Value = 12518.2
$ws.Range("N62").Value = -13766.2
$ws.Range("H65").Value = 11120868
$ws.Range("J65").Value = 12518.2
$ws.Range("L65").Value = 62591
$ws.Range("N65").Value = -68831
$ws.Range("H99").Value = 4584.9287
$ws.Range("I99").Value = 4026.7273
$ws.Range("K99").Value = 4026.7273
$ws.Range("M99").Value = -2528.7273
$ws.Range("H107").Value = 633
$ws.Range("I107").Value = 633
$ws.Range("K107").Value = 633
$ws.Range("M107").Value = 1287
$ws.Range("H126").Value = 4584.9287
$ws.Range("I126").Value = 4026.7273
$ws.Range("K126").Value = 12080.1819
$ws.Range("M126").Value = -9610.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 412350.5
$ws.Range("I9").Value = 412350.5
$ws.Range("K9").Value = 1237051.5
$ws.Range("M9").Value = -1236827.5
$ws.Range("H23").Value = 11673.286
$ws.Range("I23").Value = 40.25
$ws.Range("J23").Value = 27184
$ws.Range("K23").Value = 120.75
$ws.Range("L23").Value = 81552
$ws.Range("M23").Value = 114.25
$ws.Range("N23").Value = -82022
$ws.Range("H25").Value = 9285.286
$ws.Range("J25").Value = 9285.286
$ws.Range("L25").Value = 27855.858
$ws.Range("N25").Value = -28193.858
$ws.Range("H30").Value = 9285.286
$ws.Range("J30").Value = 9285.286
$ws.Range("L30").Value = 27855.858
$ws.Range("N30").Value = -28059.858
$ws.Range("H39").Value = 7454.8
$ws.Range("J39").Value = 8272
$ws.Range("L39").Value = 24816
$ws.Range("N39").Value = -25404
$ws.Range("H87").Value = 8126.75
$ws.Range("I87").Value = 6669
$ws.Range("K87").Value = 20007
$ws.Range("M87").Value = -18759
$ws.Range("H90").Value = 8126.75
$ws.Range("I90").Value = 6669
$ws.Range("K90").Value = 60021
$ws.Range("M90").Value = -53781
$ws.Range("H114").Value = 929.4545000000001
$ws.Range("J114").Value = 1819.8
$ws.Range("L114").Value = 5459.4
$ws.Range("N114").Value = -11967.4
$ws.Range("H116").Value = 1819332.4
$ws.Range("I116").Value = 2726499.8
$ws.Range("J116").Value = 4997.5
$ws.Range("K116").Value = 8179499.399999999
$ws.Range("L116").Value = 14992.5
$ws.Range("M116").Value = -8176057.399999999
$ws.Range("N116").Value = -21876.5
$ws.Range("H137").Value = 2638.85
$ws.Range("I137").Value = 2063.2307
$ws.Range("J137").Value = 3707.8572
$ws.Range("K137").Value = 6189.6921
$ws.Range("L137").Value = 11123.5716
$ws.Range("M137").Value = -1089.6921
$ws.Range("N137").Value = -21323.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 127.55556
$ws.Range("I2").Value = 37.666668
$ws.Range("J2").Value = 307.33334
$ws.Range("K2").Value = 37.666668
$ws.Range("L2").Value = 307.33334
$ws.Range("M2").Value = 75.333332
$ws.Range("N2").Value = -533.33334
$ws.Range("H70").Value = 122531.88
$ws.Range("I70").Value = 171090.33
$ws.Range("K70").Value = 171090.33
$ws.Range("M70").Value = -170820.33
$ws.Range("H73").Value = 122531.88
$ws.Range("I73").Value = 171090.33
$ws.Range("K73").Value = 171090.33
$ws.Range("M73").Value = -170154.33
$ws.Range("H80").Value = 90911864
$ws.Range("I80").Value = 111113550
$ws.Range("J80").Value = 4250.5
$ws.Range("K80").Value = 111113550
$ws.Range("L80").Value = 4250.5
$ws.Range("M80").Value = -111112552
$ws.Range("N80").Value = -6246.5
$ws.Range("H83").Value = 90911864
$ws.Range("I83").Value = 111113550
$ws.Range("J83").Value = 4250.5
$ws.Range("K83").Value = 555567750
$ws.Range("L83").Value = 21252.5
$ws.Range("M83").Value = -555562758
$ws.Range("N83").Value = -31236.5
$ws.Range("H99").Value = 1453.4
$ws.Range("I99").Value = 1453.4
$ws.Range("K99").Value = 1453.4
$ws.Range("M99").Value = 792.5999999999999
$ws.Range("H122").Value = 5210.591
$ws.Range("I122").Value = 3977.25
$ws.Range("K122").Value = 11931.75
$ws.Range("M122").Value = -9481.75
$ws.Range("H126").Value = 6672.909
$ws.Range("I126").Value = 2426.625
$ws.Range("K126").Value = 7279.875
$ws.Range("M126").Value = -4809.875
$ws.Range("H132").Value = 5952.4
$ws.Range("I132").Value = 1654.9166
$ws.Range("J132").Value = 12398.625
$ws.Range("K132").Value = 4964.7498
$ws.Range("L132").Value = 37195.875
$ws.Range("M132").Value = -2434.7498
$ws.Range("N132").Value = -42255.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 11999.5
$ws.Range("J17").Value = 11999.5
$ws.Range("L17").Value = 11999.5
$ws.Range("N17").Value = -12339.5
$ws.Range("H22").Value = 762.63635
$ws.Range("I22").Value = 422.75
$ws.Range("K22").Value = 422.75
$ws.Range("M22").Value = -127.75
$ws.Range("H27").Value = 762.63635
$ws.Range("I27").Value = 422.75
$ws.Range("K27").Value = 422.75
$ws.Range("M27").Value = -315.75
$ws.Range("H40").Value = 4524.7837
$ws.Range("I40").Value = 4571.909
$ws.Range("K40").Value = 4571.909
$ws.Range("M40").Value = -4435.909
$ws.Range("H55").Value = 356.08334
$ws.Range("I55").Value = 146.85715
$ws.Range("J55").Value = 649
$ws.Range("K55").Value = 146.85715
$ws.Range("L55").Value = 649
$ws.Range("M55").Value = 26.14285000000001
$ws.Range("N55").Value = -995
$ws.Range("H132").Value = 4812.4
$ws.Range("I132").Value = 2183.92
$ws.Range("J132").Value = 11383.6
$ws.Range("K132").Value = 6551.76
$ws.Range("L132").Value = 34150.8
$ws.Range("M132").Value = -4021.76
$ws.Range("N132").Value = -39210.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 887.0833
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
